$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows at 14:15 (pushes everything below down by two rows).
#    This makes room for two new "Docentes responsaveis" entries that sit
#    right under "519033 - Carlos Yujiro Shigue" (row 13) and right above
#    "Programa resumido:" (old row 14, new row 16).
# ---------------------------------------------------------------------------
$ws.Rows("14:15").Insert()

# New teacher rows - only columns B (current) and C (modified) are filled,
# column A stays empty, same as the row above them (row 13).
$ws.Range("B14").Value = "1033242 - Fábio Herbst Florenzano"
$ws.Range("C14").Value = "1033242 - Fábio Herbst Florenzano"

$ws.Range("B15").Value = "1922320 - Sebastiao Ribeiro"
$ws.Range("C15").Value = "1922320 - Sebastiao Ribeiro"

# ---------------------------------------------------------------------------
# 2. "Ativação:" date value (row 8). Both B8 and C8 hold the same text and
#    must remain plain text "01/01/2022" (not get auto-converted to a date
#    serial number), so prefix with an apostrophe to force text entry.
# ---------------------------------------------------------------------------
$ws.Range("B8").Value = "'01/01/2022"
$ws.Range("C8").Value = "'01/01/2022"

# ---------------------------------------------------------------------------
# 3. "Objetivos:" text (row 10). In the original workbook B10 and C10 both
#    point at the very same shared string, and the diff only edits that
#    shared string's text in place (no sheet1.xml cell changes) - so both
#    columns must carry the new text.
# ---------------------------------------------------------------------------
$objetivos = @"
Fornecer aos estudantes uma visão abrangente e interdisciplinar sobre materiais compósitos, além de mostrar as especificidades de cada matriz, sendo ela metálica, cerâmica ou polimérica. Ademais, deseja-se apresentar os fundamentos teóricos da mecânica de estruturas reforçadas e a partir de atividades práticas demostrar métodos de caracterização de materiais compósitos e como prepara-los.
"@
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# ---------------------------------------------------------------------------
# 4. "Programa resumido:" text (new row 16) - same shared string in B and C.
# ---------------------------------------------------------------------------
$programaResumido = @"
1.Introduçâo 2. Conceitos básicos sobre materiais compósitos, suas matrizes e seus processo de fabricação 3. Tipos de reforços 4. Compósitos nanoestruturados, naturais e híbridos 5. Mecânica da estrutura reforçada 6. Atividade prática
"@
$ws.Range("B16").Value = $programaResumido
$ws.Range("C16").Value = $programaResumido

# ---------------------------------------------------------------------------
# 5. "Programa:" text (new row 18) - same shared string in B and C.
# ---------------------------------------------------------------------------
$programa = @"
1. Conceitos básicos sobre materiais compósitos: compósitos de matriz metálica (CMM), compósitos de matriz cerâmicos (CMC) e compósitos de matriz polimérica (CMP) e nanocompósitos. 2. Tipos de Reforços: Reforços particulados, fibras curtas, fibras longas, mantas, tecidos e preformas. 3. Conceitos de Interface4. Compósitos de matriz metálica: características e processos de fabricação. 5. Compósitos de matriz cerâmica: características e processos de fabricação. 6. Compósitos de matriz polimérica: matrizes termoplásticas e termorrígidas, características físicas e químicas e processos de fabricação. 7. Compósitos nanoestruturados. 8. Compósitos Naturais. 9. Compósitos Híbridos 10. Mecânica de estruturas reforçadas. Conteúdo prático: 1. Caracterização e análise de compósitos de matriz metálica. 2. Preparação e caracterização de compósitos de matriz polimérica.(Sugestão: Considerar substituir essa parte prática pela realização do PBL descrito no item 3) 3. Visita a empresa produtora de compósitos e aulas especiais e/ou palestras com professores/pesquisadores convidados
"@
$ws.Range("B18").Value = $programa
$ws.Range("C18").Value = $programa

# ---------------------------------------------------------------------------
# 6. "Método:" text (new row 21) - same shared string in B and C.
# ---------------------------------------------------------------------------
$metodo = @"
De acordo com a atual ementa da disciplina propõe-se o uso de uma nova metodologia de ensino com o intuito de abordar o conteúdo de forma mais prática e contextualizada para que o aluno consiga relacionar os conhecimentos teóricos vistos em sala de aula com as outras disciplinas do curso. Assim, avaliação do aluno será feita através de uma prova escrita e por uma apresentação final com base nas atividades práticas desenvolvidas.
"@
$ws.Range("B21").Value = $metodo
$ws.Range("C21").Value = $metodo

# ---------------------------------------------------------------------------
# 7. "Critério:" text (new row 22) - same shared string in B and C.
# ---------------------------------------------------------------------------
$criterio = @"
A nota final será calculada como descrita a seguir: NF= (0,4*Avaliação escrita + 0,6 *Apresentação final)
"@
$ws.Range("B22").Value = $criterio
$ws.Range("C22").Value = $criterio

# ---------------------------------------------------------------------------
# 8. "Norma de recuperação:" text (new row 23) - same shared string in B/C.
# ---------------------------------------------------------------------------
$norma = @"
Devido a cunho prático da disciplina não haverá recuperação.
"@
$ws.Range("B23").Value = $norma
$ws.Range("C23").Value = $norma

# ---------------------------------------------------------------------------
# 9. "Bibliografia:" text (new row 24) - same shared string in B and C.
# ---------------------------------------------------------------------------
$biblio = @"
1. REZENDE, M. C.; COSTA, M. L.; BOTELHO, E. C. Compósitos estruturais: tecnologia e prática. São Paulo: Artliber, 2011. 396p. 2 MALLICK, P.K. Composites Engineering Handbook. New York: Marcel Dekker, 1997. 3. MATTHEWS, F.L. & RAWLINGS, R.D. Composite Materials: Engineering and Science. London: Chapman & Hall, 1994. 4. OBRAZTSOV, I.F. Mechanics of Composites. Moscow: MIR Publishers, 1982. 5. JONES R. Mechanics of Composite Materials. New York: McGraw-Hill, 1975. 6. UPADHYAYA, G.S. Sintered Metal-Ceramic Composites. Elsevier, 1984. 7. HARPER, C. A. Handbook of Plastics, Elastomers and Composites. New York: McGraw-Hill, 1992. 8. GOLDSTEIN, A.N. Handbook of Nanophase Materials. CRC Press, 1997. 9. DRESSELHAUS, M.S. Graphite Fibers and Filaments. New York: Springer-Verlag, 1988.
"@
$ws.Range("B24").Value = $biblio
$ws.Range("C24").Value = $biblio

Write-Host "Edits applied."
